$wb = $excel.ActiveWorkbook

# The "Set" sheet had its "timestamp" column (column B) removed, which
# shifts "comedian" and "venue_id" one column to the left (B and C).
$setSheet = $wb.Worksheets.Item("Set")
$setSheet.Columns("B:B").Delete()

# After the edit, the author's selection on the Set sheet sits at D13.
$setSheet.Range("D13").Select()

# The active/selected tab moved from "SwearWord" to "Venue".
$venueSheet = $wb.Worksheets.Item("Venue")
$venueSheet.Activate()
